$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# segment-name column (A) to B, PercActivations (B) to C, etc.
$ws.Columns("A:A").Insert()

# The former segment-name column (now column B) already carries the
# header-ish style (s="1"); copy that same formatting onto the new
# index column A2:A20.
$ws.Range("B2").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)

# Give the new "segments" header cell (B1) the same style as the other
# header cells (C1:F1) before writing its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# Fill the new column A with a 0-based row index for each data row.
for ($i = 0; $i -lt 19; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}
